# Scheduled market-data refresh for the Typhon_Profits workbook.
# Updates per-Leve price/profit figures (columns H-N) on each Disciple
# of the Hand sheet, pulled from the latest Universalis market snapshot.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 1542.1052
$ws.Range("J69").Value = 1500
$ws.Range("L69").Value = 4500
$ws.Range("N69").Value = -6248
$ws.Range("H72").Value = 1542.1052
$ws.Range("J72").Value = 1500
$ws.Range("L72").Value = 13500
$ws.Range("N72").Value = -22236
$ws.Range("H116").Value = 4813.7144
$ws.Range("I116").Value = 2300
$ws.Range("J116").Value = 5232.6665
$ws.Range("K116").Value = 2300
$ws.Range("L116").Value = 5232.6665
$ws.Range("M116").Value = 1142
$ws.Range("N116").Value = -12116.6665
$ws.Range("H129").Value = 264132.5
$ws.Range("J129").Value = 346036.44
$ws.Range("L129").Value = 1038109.32
$ws.Range("N129").Value = -1048109.32
$ws.Range("H132").Value = 1965.6271
$ws.Range("I132").Value = 2060.7693
$ws.Range("J132").Value = 1258.8572
$ws.Range("K132").Value = 6182.3079
$ws.Range("L132").Value = 3776.5716
$ws.Range("M132").Value = -3652.3079
$ws.Range("N132").Value = -8836.571599999999
$ws.Range("H137").Value = 1276.9259
$ws.Range("I137").Value = 1288.85
$ws.Range("J137").Value = 1242.8572
$ws.Range("K137").Value = 3866.55
$ws.Range("L137").Value = 3728.5716
$ws.Range("M137").Value = -1316.55
$ws.Range("N137").Value = -8828.571599999999
$ws.Range("H138").Value = 25643318
$ws.Range("I138").Value = 52632904
$ws.Range("J138").Value = 3210.4
$ws.Range("K138").Value = 157898712
$ws.Range("L138").Value = 9631.200000000001
$ws.Range("M138").Value = -157893572
$ws.Range("N138").Value = -19911.2
$ws.Range("H141").Value = 3075.889
$ws.Range("I141").Value = 2613.8333
$ws.Range("J141").Value = 4000
$ws.Range("K141").Value = 7841.499899999999
$ws.Range("L141").Value = 12000
$ws.Range("M141").Value = -2661.499899999999
$ws.Range("N141").Value = -22360

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9949.159
$ws.Range("I32").Value = 7161.4
$ws.Range("J32").Value = 20790.445
$ws.Range("K32").Value = 7161.4
$ws.Range("L32").Value = 20790.445
$ws.Range("M32").Value = -6874.4
$ws.Range("N32").Value = -21364.445
$ws.Range("H61").Value = 2402.5
$ws.Range("I61").Value = 1550
$ws.Range("K61").Value = 1550
$ws.Range("M61").Value = -1338
$ws.Range("H74").Value = 25001278
$ws.Range("I74").Value = 33333718
$ws.Range("K74").Value = 33333718
$ws.Range("M74").Value = -33332844
$ws.Range("H77").Value = 25001278
$ws.Range("I77").Value = 33333718
$ws.Range("K77").Value = 166668590
$ws.Range("M77").Value = -166664222
$ws.Range("H109").Value = 20000
$ws.Range("J109").Value = 20000
$ws.Range("L109").Value = 20000
$ws.Range("N109").Value = -22774
$ws.Range("H114").Value = 29579.4
$ws.Range("J114").Value = 29579.4
$ws.Range("L114").Value = 29579.4
$ws.Range("N114").Value = -38257.4
$ws.Range("H119").Value = 31232.666
$ws.Range("J119").Value = 31232.666
$ws.Range("L119").Value = 31232.666
$ws.Range("N119").Value = -40908.666
$ws.Range("H132").Value = 25703.857
$ws.Range("I132").Value = 2121.7693
$ws.Range("J132").Value = 64024.75
$ws.Range("K132").Value = 6365.3079
$ws.Range("L132").Value = 192074.25
$ws.Range("M132").Value = -3835.3079
$ws.Range("N132").Value = -197134.25
$ws.Range("H136").Value = 2402.5
$ws.Range("I136").Value = 1550
$ws.Range("K136").Value = 4650
$ws.Range("M136").Value = -2100

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 1749.5
$ws.Range("I11").Value = 1749.5
$ws.Range("K11").Value = 1749.5
$ws.Range("M11").Value = -1609.5
$ws.Range("H110").Value = 31850
$ws.Range("J110").Value = 31850
$ws.Range("L110").Value = 31850
$ws.Range("N110").Value = -40030
$ws.Range("H138").Value = 50676
$ws.Range("J138").Value = 50676
$ws.Range("L138").Value = 50676
$ws.Range("N138").Value = -60956

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 15000
$ws.Range("J43").Value = 15000
$ws.Range("L43").Value = 15000
$ws.Range("N43").Value = -15368
$ws.Range("H58").Value = 20154.74
$ws.Range("I58").Value = 1713.8334
$ws.Range("J58").Value = 34907.465
$ws.Range("K58").Value = 1713.8334
$ws.Range("L58").Value = 34907.465
$ws.Range("M58").Value = -1510.8334
$ws.Range("N58").Value = -35313.465
$ws.Range("H95").Value = 20000
$ws.Range("J95").Value = 20000
$ws.Range("L95").Value = 20000
$ws.Range("N95").Value = -25492
$ws.Range("H101").Value = 15000
$ws.Range("J101").Value = 15000
$ws.Range("L101").Value = 15000
$ws.Range("N101").Value = -21490
$ws.Range("H132").Value = 2543.7778
$ws.Range("I132").Value = 1987.3684
$ws.Range("J132").Value = 3865.25
$ws.Range("K132").Value = 5962.1052
$ws.Range("L132").Value = 11595.75
$ws.Range("M132").Value = -3432.1052
$ws.Range("N132").Value = -16655.75
$ws.Range("H136").Value = 20154.74
$ws.Range("I136").Value = 1713.8334
$ws.Range("J136").Value = 34907.465
$ws.Range("K136").Value = 5141.5002
$ws.Range("L136").Value = 104722.395
$ws.Range("M136").Value = -2591.5002
$ws.Range("N136").Value = -109822.395

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 713.64
$ws.Range("J131").Value = 717.8182
$ws.Range("L131").Value = 2153.4546
$ws.Range("N131").Value = -12233.4546

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4254.375
$ws.Range("J80").Value = 4825
$ws.Range("L80").Value = 4825
$ws.Range("N80").Value = -6821
$ws.Range("H83").Value = 4254.375
$ws.Range("J83").Value = 4825
$ws.Range("L83").Value = 24125
$ws.Range("N83").Value = -34109
$ws.Range("H102").Value = 15153002
$ws.Range("I102").Value = 17242580
$ws.Range("K102").Value = 17242580
$ws.Range("M102").Value = -17240958
$ws.Range("H107").Value = 2849248.2
$ws.Range("I107").Value = 258.70587
$ws.Range("J107").Value = 7692530.5
$ws.Range("K107").Value = 258.70587
$ws.Range("L107").Value = 7692530.5
$ws.Range("M107").Value = 1661.29413
$ws.Range("N107").Value = -7696370.5
$ws.Range("H109").Value = 28875
$ws.Range("J109").Value = 28875
$ws.Range("L109").Value = 28875
$ws.Range("N109").Value = -30955
$ws.Range("H126").Value = 3916.111
$ws.Range("I126").Value = 2607.5
$ws.Range("K126").Value = 7822.5
$ws.Range("M126").Value = -5352.5
$ws.Range("H132").Value = 36219.133
$ws.Range("I132").Value = 2602.375
$ws.Range("K132").Value = 7807.125
$ws.Range("M132").Value = -5277.125
$ws.Range("H140").Value = 40000
$ws.Range("J140").Value = 40000
$ws.Range("L140").Value = 40000
$ws.Range("N140").Value = -50360

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3944.4285
$ws.Range("I61").Value = 1774.7273
$ws.Range("K61").Value = 1774.7273
$ws.Range("M61").Value = -1572.7273
$ws.Range("H113").Value = 3944.4285
$ws.Range("I113").Value = 1774.7273
$ws.Range("K113").Value = 1774.7273
$ws.Range("M113").Value = 395.2727
$ws.Range("H136").Value = 1337.7142
$ws.Range("I136").Value = 1268
$ws.Range("K136").Value = 3804
$ws.Range("M136").Value = -1254

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H132").Value = 1127.7826
$ws.Range("I132").Value = 908.3158
$ws.Range("K132").Value = 2724.9474
$ws.Range("M132").Value = -194.9474

